# Trading update: 2026-02-18 10:52:11
#
# - The "MarketMaking" strategy's existing trade-tracker sheet (with its
#   closed trade #1 and open trade #3) is retired/dropped.
# - The "momentum" strategy sheet -- which only ever held trade #2 -- gets
#   taken over by MarketMaking: it's renamed "MarketMaking" and its single
#   row is overwritten with the strategy's brand-new trade (#4).
# - That new trade is also appended to "All Trades" as row 5.
# - Older rows on "All Trades" lose their live analytics columns (Capital
#   After, slippage, confidence, entry/exit reason, duration) now that
#   they're no longer the newest trade for their strategy; OPEN trades
#   without an exit yet show a numeric 0 Exit Price instead of being blank.
# - "Strategy Status" no longer counts the trade that rolled off
#   MarketMaking's old tracker sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Drop the old "MarketMaking" sheet, then repurpose "momentum" in its
#    place (same slot the old sheet occupied, same sheet tab name).
# ---------------------------------------------------------------------------
$oldMarketMaking = $wb.Worksheets.Item("MarketMaking")
$oldMarketMaking.Delete()

$mm = $wb.Worksheets.Item("momentum")
$mm.Name = "MarketMaking"

# ---------------------------------------------------------------------------
# 2. Strategy Status: MarketMaking's trade counter resets to 0.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 0

# ---------------------------------------------------------------------------
# 3. MarketMaking sheet: overwrite its single open-trade row with the new
#    trade (#4) that just opened.
# ---------------------------------------------------------------------------
$mm.Range("A2").Value = 4
$mm.Range("C2").Value = "10:50:40"
$mm.Range("D2").Value = "MarketMaking"
$mm.Range("E2").Value = "DOWN"
$mm.Range("F2").Value = 0.45
$mm.Range("G2").ClearContents()
$mm.Range("H2").Value = "OPEN"
$mm.Range("I2").Value = 0
$mm.Range("J2").Value = 0
$mm.Range("K2").Value = 100
$mm.Range("L2").Value = 0
$mm.Range("M2").Value = 0
$mm.Range("N2").Value = 0.6
$mm.Range("O2").Value = "Normal spread capture: 202 bps"
$mm.Range("P2").ClearContents()
$mm.Range("Q2").Value = 0

# ---------------------------------------------------------------------------
# 4. All Trades: older rows lose their per-trade analytics now that they are
#    no longer the latest trade for their strategy, then append the new
#    MarketMaking trade as row 5.
# ---------------------------------------------------------------------------
$all = $wb.Worksheets.Item("All Trades")

# Row 2 - Trade #1 (MarketMaking, CLOSED): clear analytics columns.
$all.Range("K2:Q2").ClearContents()

# Row 3 - Trade #2 (momentum, OPEN): Exit Price becomes an explicit 0,
# analytics columns clear out.
$all.Range("G3").Value = 0
$all.Range("K3:O3").ClearContents()
$all.Range("Q3").ClearContents()

# Row 4 - Trade #3 (MarketMaking, OPEN): same treatment as row 3.
$all.Range("G4").Value = 0
$all.Range("K4:O4").ClearContents()
$all.Range("Q4").ClearContents()

# Row 5 - Trade #4 (MarketMaking, OPEN): newly appended trade, full detail.
# (Date/time are stored as plain text, not real Excel dates/times -- the
# leading apostrophe forces text entry the same way typing it in the UI
# would, matching the rest of the sheet's Date/Time columns.)
$all.Range("A5").Value = 4
$all.Range("B5").Value = "'2026-02-18"
$all.Range("C5").Value = "'10:50:40"
$all.Range("D5").Value = "MarketMaking"
$all.Range("E5").Value = "DOWN"
$all.Range("F5").Value = 0.45
$all.Range("G5").ClearContents()
$all.Range("H5").Value = "OPEN"
$all.Range("I5").Value = 0
$all.Range("J5").Value = 0
$all.Range("K5").Value = 100
$all.Range("L5").Value = 0
$all.Range("M5").Value = 0
$all.Range("N5").Value = 0.6
$all.Range("O5").Value = "Normal spread capture: 202 bps"
$all.Range("P5").ClearContents()
$all.Range("Q5").Value = 0

# ---------------------------------------------------------------------------
# 5. Restore the original active sheet/selection (deleting a sheet shifts
#    Excel's active-tab focus as a side effect).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Summary").Activate()
